# Insert a new data row before row 371 (Feria Lagunitas de Puerto Montt - Ajo),
# shifting the existing rows 371..499 down to 372..500, then populate the new
# row 371 with its own record (D/K/L/M/P differ from the old row 371 that is
# now at row 372; the rest of the fields are carried over unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("371:371").Insert()

$ws.Range("A371").Value2 = 4
$ws.Range("B371").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C371").Value2 = "Los Lagos"
$ws.Range("D371").Value2 = 45146
$ws.Range("E371").Value2 = 10
$ws.Range("F371").Value2 = 100112003
$ws.Range("G371").Value2 = "Ajo"
$ws.Range("H371").Value2 = "Chino"
$ws.Range("I371").Value2 = "Primera"
$ws.Range("J371").Value2 = 220
$ws.Range("K371").Value2 = 23000
$ws.Range("L371").Value2 = 23000
$ws.Range("M371").Value2 = 23000
$ws.Range("N371").Value2 = "`$/caja 10 kilos"
$ws.Range("O371").Value2 = "China"
$ws.Range("P371").Value2 = 2300
$ws.Range("Q371").Value2 = 10
$ws.Range("R371").Value2 = "Hortaliza"
